$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 646
$ws1.Range("F9").Value = 403
$ws1.Range("F10").Value = 4065
$ws1.Range("F12").Value = 437
$ws1.Range("F18").Value = 2924
$ws1.Range("F24").Value = 902
$ws1.Range("F26").Value = 2174
$ws1.Range("F27").Value = 967
$ws1.Range("F28").Value = 2260
$ws1.Range("F31").Value = 428
$ws1.Range("F34").Value = 386
$ws1.Range("F35").Value = 1034
$ws1.Range("F37").Value = 1128
$ws1.Range("F38").Value = 286

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 646
$ws4.Range("F10").Value = 403
$ws4.Range("F11").Value = 4065
$ws4.Range("F19").Value = 2924
$ws4.Range("F29").Value = 902
$ws4.Range("F30").Value = 2174
$ws4.Range("F33").Value = 967
$ws4.Range("F34").Value = 2260
$ws4.Range("F36").Value = 428
$ws4.Range("F38").Value = 386
$ws4.Range("F39").Value = 1034
$ws4.Range("F41").Value = 1128
$ws4.Range("F42").Value = 286
